# "Prop book.xlsx" update: add an "AFS value" column, reclassify VCI's
# existing 2Q25/1Q25/4Q24 figures from FVTPL to AFS, and append new VDS
# and DSE broker rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn off iterative calculation (workbook-level setting)
$excel.Iteration = $false

# 1. Append the new VDS / DSE broker rows (209-238) first, so the new
#    ticker strings land in the shared-string table in this order:
#    VDS, KBC, CMG, QNS, DSE.
$ws.Range("A209").Value = "VDS"
$ws.Range("B209").Value = "2Q25"
$ws.Range("C209").Value = "KBC"
$ws.Range("D209").Value = 291.53755
$ws.Range("A210").Value = "VDS"
$ws.Range("B210").Value = "2Q25"
$ws.Range("C210").Value = "HSG"
$ws.Range("D210").Value = 97.5997
$ws.Range("A211").Value = "VDS"
$ws.Range("B211").Value = "2Q25"
$ws.Range("C211").Value = "MWG"
$ws.Range("D211").Value = 139.9604
$ws.Range("A212").Value = "VDS"
$ws.Range("B212").Value = "2Q25"
$ws.Range("C212").Value = "ACB"
$ws.Range("D212").Value = 126.0891627
$ws.Range("A213").Value = "VDS"
$ws.Range("B213").Value = "2Q25"
$ws.Range("C213").Value = "Others"
$ws.Range("D213").Value = 484.16045325
$ws.Range("A214").Value = "VDS"
$ws.Range("B214").Value = "2Q25"
$ws.Range("C214").Value = "CMG"
$ws.Range("E214").Value = 97.86035
$ws.Range("A215").Value = "VDS"
$ws.Range("B215").Value = "2Q25"
$ws.Range("C215").Value = "KDH"
$ws.Range("E215").Value = 118.3644
$ws.Range("A216").Value = "VDS"
$ws.Range("B216").Value = "2Q25"
$ws.Range("C216").Value = "KBC"
$ws.Range("E216").Value = 40.125
$ws.Range("A217").Value = "VDS"
$ws.Range("B217").Value = "2Q25"
$ws.Range("C217").Value = "PBT"
$ws.Range("D217").Value = -17.442687286
$ws.Range("A218").Value = "VDS"
$ws.Range("B218").Value = "1Q25"
$ws.Range("C218").Value = "KBC"
$ws.Range("D218").Value = 319.48334
$ws.Range("A219").Value = "VDS"
$ws.Range("B219").Value = "1Q25"
$ws.Range("C219").Value = "Others"
$ws.Range("D219").Value = 1110.0100874
$ws.Range("A220").Value = "VDS"
$ws.Range("B220").Value = "1Q25"
$ws.Range("C220").Value = "KDH"
$ws.Range("E220").Value = 131.2476
$ws.Range("A221").Value = "VDS"
$ws.Range("B221").Value = "1Q25"
$ws.Range("C221").Value = "CMG"
$ws.Range("E221").Value = 77.8012
$ws.Range("A222").Value = "VDS"
$ws.Range("B222").Value = "1Q25"
$ws.Range("C222").Value = "QNS"
$ws.Range("E222").Value = 49.611
$ws.Range("A223").Value = "VDS"
$ws.Range("B223").Value = "1Q25"
$ws.Range("C223").Value = "MWG"
$ws.Range("E223").Value = 36.58
$ws.Range("A224").Value = "VDS"
$ws.Range("B224").Value = "1Q25"
$ws.Range("C224").Value = "PBT"
$ws.Range("D224").Value = 22.643851018
$ws.Range("A225").Value = "VDS"
$ws.Range("B225").Value = "4Q24"
$ws.Range("C225").Value = "VNM"
$ws.Range("D225").Value = 147.13238
$ws.Range("A226").Value = "VDS"
$ws.Range("B226").Value = "4Q24"
$ws.Range("C226").Value = "HSG"
$ws.Range("D226").Value = 102.867975
$ws.Range("A227").Value = "VDS"
$ws.Range("B227").Value = "4Q24"
$ws.Range("C227").Value = "CTG"
$ws.Range("D227").Value = 125.2624716
$ws.Range("A228").Value = "VDS"
$ws.Range("B228").Value = "4Q24"
$ws.Range("C228").Value = "ACB"
$ws.Range("D228").Value = 129.7107126
$ws.Range("A229").Value = "VDS"
$ws.Range("B229").Value = "4Q24"
$ws.Range("C229").Value = "MWG"
$ws.Range("D229").Value = 127.5998
$ws.Range("A230").Value = "VDS"
$ws.Range("B230").Value = "4Q24"
$ws.Range("C230").Value = "Others"
$ws.Range("D230").Value = 466.7749656
$ws.Range("A231").Value = "VDS"
$ws.Range("B231").Value = "4Q24"
$ws.Range("C231").Value = "KBC"
$ws.Range("E231").Value = 270.7216
$ws.Range("A232").Value = "VDS"
$ws.Range("B232").Value = "4Q24"
$ws.Range("C232").Value = "KDH"
$ws.Range("E232").Value = 145.3386
$ws.Range("A233").Value = "VDS"
$ws.Range("B233").Value = "4Q24"
$ws.Range("C233").Value = "QNS"
$ws.Range("E233").Value = 50.883
$ws.Range("A234").Value = "VDS"
$ws.Range("B234").Value = "4Q24"
$ws.Range("C234").Value = "Others"
$ws.Range("E234").Value = 62.96705
$ws.Range("A235").Value = "VDS"
$ws.Range("B235").Value = "4Q24"
$ws.Range("C235").Value = "PBT"
$ws.Range("D235").Value = -33.63417417
$ws.Range("A236").Value = "DSE"
$ws.Range("B236").Value = "2Q25"
$ws.Range("C236").Value = "PBT"
$ws.Range("D236").Value = 91.021986688
$ws.Range("A237").Value = "DSE"
$ws.Range("B237").Value = "1Q25"
$ws.Range("C237").Value = "PBT"
$ws.Range("D237").Value = 66.466997524
$ws.Range("A238").Value = "DSE"
$ws.Range("B238").Value = "4Q24"
$ws.Range("C238").Value = "PBT"
$ws.Range("D238").Value = 41.544733349

# 2. Add the new "AFS value" column header
$ws.Range("E1").Value = "AFS value"

# 3. Re-home the VCI broker block (rows 2-28, all three quarters) from the
#    FVTPL value column (D) into the new AFS value column (E).
$ws.Range("E2").Value = 598.5
$ws.Range("D2").Clear()
$ws.Range("E3").Value = 81.7
$ws.Range("D3").Clear()
$ws.Range("E4").Value = 1945.2
$ws.Range("D4").Clear()
$ws.Range("E6").Value = 180.69
$ws.Range("D6").Clear()
$ws.Range("E7").Value = 3.3
$ws.Range("D7").Clear()
$ws.Range("E8").Value = 232.999
$ws.Range("D8").Clear()
$ws.Range("E9").Value = 853
$ws.Range("D9").Clear()
$ws.Range("E10").Value = 2756
$ws.Range("D10").Clear()
$ws.Range("E11").Value = 796.99
$ws.Range("D11").Clear()
$ws.Range("E13").Value = 1882.8
$ws.Range("D13").Clear()
$ws.Range("E14").Value = 11.73
$ws.Range("D14").Clear()
$ws.Range("E15").Value = 297.25
$ws.Range("D15").Clear()
$ws.Range("E16").Value = 139.4
$ws.Range("D16").Clear()
$ws.Range("E17").Value = 803.12
$ws.Range("D17").Clear()
$ws.Range("E18").Value = 837.54
$ws.Range("D18").Clear()
$ws.Range("E19").Value = 1673.2
$ws.Range("D19").Clear()
$ws.Range("E20").Value = 1184
$ws.Range("D20").Clear()
$ws.Range("E22").Value = 2117.6
$ws.Range("D22").Clear()
$ws.Range("E24").Value = 12.33
$ws.Range("D24").Clear()
$ws.Range("E25").Value = 167.5
$ws.Range("D25").Clear()
$ws.Range("E26").Value = 695.36
$ws.Range("D26").Clear()
$ws.Range("E27").Value = 772.4
$ws.Range("D27").Clear()
$ws.Range("E28").Value = 2070
$ws.Range("D28").Clear()

# 4. Match the saved cursor position recorded in the workbook
$ws.Range("B3").Select() | Out-Null
